# Fixed data - LDA
# The sign of the LD1 component (column A) was flipped for every data row.
# Column B (LD2) is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = -1 * $val
    }
}
